$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: extend sequence with P1=14, Q1=15, matching style of O1 ---
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows 2-25: swap I/K and M/O columns, then add new P/Q columns ---
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value()   # column I
    $kVal = $ws.Cells.Item($r, 11).Value()  # column K
    $mVal = $ws.Cells.Item($r, 13).Value()  # column M
    $oVal = $ws.Cells.Item($r, 15).Value()  # column O

    $ws.Cells.Item($r, 9).Value = $kVal    # I <- old K
    $ws.Cells.Item($r, 11).Value = $iVal   # K <- old I
    $ws.Cells.Item($r, 13).Value = $oVal   # M <- old O
    $ws.Cells.Item($r, 15).Value = $mVal   # O <- old M

    $ws.Cells.Item($r, 16).Value = 2       # P
    $ws.Cells.Item($r, 17).Value = 2       # Q
}
